$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "online"
$ws.Range("F2").Value = $false

$ws.Range("F2").Select()
